# Add season record (Wins/Losses/Ties) columns to the right of the
# existing team stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, bordered, centered) used by the
# rest of row 1 by copying the style from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Season record values for every player row (2-50): 70 wins, 92 losses,
# 0 ties.
$lastRow = 50
$wins = $ws.Range("AD2:AD" + $lastRow)
$losses = $ws.Range("AE2:AE" + $lastRow)
$ties = $ws.Range("AF2:AF" + $lastRow)

$wins.Value = 70
$losses.Value = 92
$ties.Value = 0
